$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" date field text on the slide
#    master and on every slide layout (the template's footer/date
#    placeholder) from "6/5/17" to "2/17/2018".
# ---------------------------------------------------------------------------
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $ph = $shp.PlaceholderFormat
            if ($ph -ne $null -and $ph.Type -eq $ppPlaceholderDate) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq "6/5/17") {
                    $tr.Text = "2/17/2018"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Add the two presentation-level slide guides (a horizontal guide at
#    2160 and a vertical guide at 2880) that PowerPoint records in
#    p:presentation/p:extLst as a p15:sldGuideLst.
# ---------------------------------------------------------------------------
$ppHorizontalGuide = 1
$ppVerticalGuide = 2

$guides = $p.Guides
$guides.Add($ppHorizontalGuide, 2160)
$guides.Add($ppVerticalGuide, 2880)
